$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells keep their exact text representation
# (Excel would otherwise parse plain-looking numeric strings into floating
# point numbers and lose trailing zeros / exact formatting), so force the
# number format to Text before writing the new values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "26.724.31"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.601.56"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "211.60"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "0.247"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").Value = "1.826.92"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.601.44"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "64.93"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "26.695.63"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "0.0₃0742"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "210.25"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "7.18"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("D24").Value = "8.97"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "144.21"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "7.09"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("D29").Value = "15.37"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "0.0511"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Value = "1.295.37"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("E37").Value = "  -2.62%  "
$ws.Range("D38").Value = "1.15"
$ws.Range("E38").Value = "  +7.84%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "0.829"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").Value = "5.40"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "2.19"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "0.781"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Value = "63.04"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "1.739.13"
$ws.Range("D46").Value = "90.77"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("E47").Value = "  -1.95%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("E51").Value = "  -0.26%  "
